$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (232-240) transcribing the additional verse analysis
# "ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥" word-by-word, one row per word,
# matching the existing sheet layout (columns A-AC).

# Row 232
$ws.Cells.Item(232, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(232, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(232, 3).Value = ''
$ws.Cells.Item(232, 4).Value = 'ਕਬੀਰ'
$ws.Cells.Item(232, 5).Value = 'ਹੇ ਕਬੀਰ!'
$ws.Cells.Item(232, 6).Value = 'ਮੁਕਤਾ'
$ws.Cells.Item(232, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(232, 8).Value = 'Vocative ਸੰਬੋਧਨ'
$ws.Cells.Item(232, 9).Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Cells.Item(232, 10).Value = 'ਮੁਕਤਾ Ending'
$ws.Cells.Item(232, 11).Value = 'Noun / ਨਾਂਵ'
$ws.Cells.Item(232, 12).Value = 1
$ws.Cells.Item(232, 13).Value = 0
$ws.Cells.Item(232, 14).Value = 54730
$ws.Cells.Item(232, 15).Value = 1
$ws.Cells.Item(232, 16).Value = 1
$ws.Cells.Item(232, 17).Value = 3818
$ws.Cells.Item(232, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(232, 19).Value = ''
$ws.Cells.Item(232, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(232, 21).Value = ''
$ws.Cells.Item(232, 22).Value = ''
$ws.Cells.Item(232, 23).Value = ''
$ws.Cells.Item(232, 24).Value = ''
$ws.Cells.Item(232, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(232, 26).Value = 'Simple'
$ws.Cells.Item(232, 27).Value = 1376
$ws.Cells.Item(232, 28).Value = 0
$ws.Cells.Item(232, 29).Value = 0

# Row 233
$ws.Cells.Item(233, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(233, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(233, 3).Value = ''
$ws.Cells.Item(233, 4).Value = 'ਮਨੁ'
$ws.Cells.Item(233, 5).Value = 'ਮਨ'
$ws.Cells.Item(233, 6).Value = 'ਮਨੁ'
$ws.Cells.Item(233, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(233, 8).Value = 'Nominative ਕਰਤਾ'
$ws.Cells.Item(233, 9).Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Cells.Item(233, 10).Value = 'ਮੁਕਤਾ Ending'
$ws.Cells.Item(233, 11).Value = 'Noun / ਨਾਂਵ'
$ws.Cells.Item(233, 12).Value = 1
$ws.Cells.Item(233, 13).Value = 1
$ws.Cells.Item(233, 14).Value = 54730
$ws.Cells.Item(233, 15).Value = 1
$ws.Cells.Item(233, 16).Value = 1
$ws.Cells.Item(233, 17).Value = 3818
$ws.Cells.Item(233, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(233, 19).Value = ''
$ws.Cells.Item(233, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(233, 21).Value = ''
$ws.Cells.Item(233, 22).Value = ''
$ws.Cells.Item(233, 23).Value = ''
$ws.Cells.Item(233, 24).Value = ''
$ws.Cells.Item(233, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(233, 26).Value = 'Simple'
$ws.Cells.Item(233, 27).Value = 1376
$ws.Cells.Item(233, 28).Value = 0
$ws.Cells.Item(233, 29).Value = 0

# Row 234
$ws.Cells.Item(234, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(234, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(234, 3).Value = ''
$ws.Cells.Item(234, 4).Value = 'ਜਾਨੈ'
$ws.Cells.Item(234, 5).Value = 'ਜਾਣਦਾ ਹੈ| ਸਮਝਦਾ ਹੈ| ਪਛਾਣਦਾ ਹੈ'
$ws.Cells.Item(234, 6).Value = 'ੈ'
$ws.Cells.Item(234, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(234, 8).Value = 'Present ਵਰਤਮਾਨ'
$ws.Cells.Item(234, 9).Value = 'Trans / ਨਪੁਂਸਕ'
$ws.Cells.Item(234, 10).Value = '3rd Person / ਅਨਯ ਪੁਰਖ'
$ws.Cells.Item(234, 11).Value = 'Verb / ਕਿਰਿਆ'
$ws.Cells.Item(234, 12).Value = 1
$ws.Cells.Item(234, 13).Value = 2
$ws.Cells.Item(234, 14).Value = 54730
$ws.Cells.Item(234, 15).Value = 1
$ws.Cells.Item(234, 16).Value = 1
$ws.Cells.Item(234, 17).Value = 3818
$ws.Cells.Item(234, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(234, 19).Value = ''
$ws.Cells.Item(234, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(234, 21).Value = ''
$ws.Cells.Item(234, 22).Value = ''
$ws.Cells.Item(234, 23).Value = ''
$ws.Cells.Item(234, 24).Value = ''
$ws.Cells.Item(234, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(234, 26).Value = 'Simple'
$ws.Cells.Item(234, 27).Value = 1376
$ws.Cells.Item(234, 28).Value = 0
$ws.Cells.Item(234, 29).Value = 0

# Row 235
$ws.Cells.Item(235, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(235, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(235, 3).Value = ''
$ws.Cells.Item(235, 4).Value = 'ਸਭ'
$ws.Cells.Item(235, 5).Value = 'ਹਰੇਕ| ਸਾਰੀ'
$ws.Cells.Item(235, 6).Value = 'ਸਭ'
$ws.Cells.Item(235, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(235, 8).Value = 'Indefinite / ਅਨਿਸਚੇ ਵਾਚਕ'
$ws.Cells.Item(235, 9).Value = 'Feminine / ਇਸਤਰੀ'
$ws.Cells.Item(235, 10).Value = ''
$ws.Cells.Item(235, 11).Value = 'Adjectives / ਵਿਸ਼ੇਸ਼ਣ'
$ws.Cells.Item(235, 12).Value = 1
$ws.Cells.Item(235, 13).Value = 3
$ws.Cells.Item(235, 14).Value = 54730
$ws.Cells.Item(235, 15).Value = 1
$ws.Cells.Item(235, 16).Value = 1
$ws.Cells.Item(235, 17).Value = 3818
$ws.Cells.Item(235, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(235, 19).Value = ''
$ws.Cells.Item(235, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(235, 21).Value = ''
$ws.Cells.Item(235, 22).Value = ''
$ws.Cells.Item(235, 23).Value = ''
$ws.Cells.Item(235, 24).Value = ''
$ws.Cells.Item(235, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(235, 26).Value = 'Simple'
$ws.Cells.Item(235, 27).Value = 1376
$ws.Cells.Item(235, 28).Value = 0
$ws.Cells.Item(235, 29).Value = 0

# Row 236
$ws.Cells.Item(236, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(236, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(236, 3).Value = ''
$ws.Cells.Item(236, 4).Value = 'ਬਾਤ'
$ws.Cells.Item(236, 5).Value = 'ਗੱਲ| ਕਹਾਣੀ'
$ws.Cells.Item(236, 6).Value = 'ਮੁਕਤਾ'
$ws.Cells.Item(236, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(236, 8).Value = 'Accusative ਕਰਮ'
$ws.Cells.Item(236, 9).Value = 'Feminine / ਇਸਤਰੀ'
$ws.Cells.Item(236, 10).Value = 'ਮੁਕਤਾ Ending'
$ws.Cells.Item(236, 11).Value = 'Noun / ਨਾਂਵ'
$ws.Cells.Item(236, 12).Value = 1
$ws.Cells.Item(236, 13).Value = 4
$ws.Cells.Item(236, 14).Value = 54730
$ws.Cells.Item(236, 15).Value = 1
$ws.Cells.Item(236, 16).Value = 1
$ws.Cells.Item(236, 17).Value = 3818
$ws.Cells.Item(236, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(236, 19).Value = ''
$ws.Cells.Item(236, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(236, 21).Value = ''
$ws.Cells.Item(236, 22).Value = ''
$ws.Cells.Item(236, 23).Value = ''
$ws.Cells.Item(236, 24).Value = ''
$ws.Cells.Item(236, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(236, 26).Value = 'Simple'
$ws.Cells.Item(236, 27).Value = 1376
$ws.Cells.Item(236, 28).Value = 0
$ws.Cells.Item(236, 29).Value = 0

# Row 237
$ws.Cells.Item(237, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(237, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(237, 3).Value = ''
$ws.Cells.Item(237, 4).Value = 'ਜਾਨਤ'
$ws.Cells.Item(237, 5).Value = 'ਜਾਣਦਾ ਹੋਇਆ'
$ws.Cells.Item(237, 6).Value = 'ਕਰਤ'
$ws.Cells.Item(237, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(237, 8).Value = 'Present ਵਰਤਮਾਨ (Participle form)'
$ws.Cells.Item(237, 9).Value = 'Trans / ਨਪੁਂਸਕ'
$ws.Cells.Item(237, 10).Value = '3rd Person / ਅਨਯ ਪੁਰਖ'
$ws.Cells.Item(237, 11).Value = 'Verb / ਕਿਰਿਆ'
$ws.Cells.Item(237, 12).Value = 1
$ws.Cells.Item(237, 13).Value = 5
$ws.Cells.Item(237, 14).Value = 54730
$ws.Cells.Item(237, 15).Value = 1
$ws.Cells.Item(237, 16).Value = 1
$ws.Cells.Item(237, 17).Value = 3818
$ws.Cells.Item(237, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(237, 19).Value = ''
$ws.Cells.Item(237, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(237, 21).Value = ''
$ws.Cells.Item(237, 22).Value = ''
$ws.Cells.Item(237, 23).Value = ''
$ws.Cells.Item(237, 24).Value = ''
$ws.Cells.Item(237, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(237, 26).Value = 'Simple'
$ws.Cells.Item(237, 27).Value = 1376
$ws.Cells.Item(237, 28).Value = 0
$ws.Cells.Item(237, 29).Value = 0

# Row 238
$ws.Cells.Item(238, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(238, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(238, 3).Value = ''
$ws.Cells.Item(238, 4).Value = 'ਹੀ'
$ws.Cells.Item(238, 5).Value = 'ਹੀ'
$ws.Cells.Item(238, 6).Value = 'ਹੀ'
$ws.Cells.Item(238, 7).Value = ''
$ws.Cells.Item(238, 8).Value = 'Sress / Emphasis / ਤਾਕੀਦ ਵਾਚਕ'
$ws.Cells.Item(238, 9).Value = ''
$ws.Cells.Item(238, 10).Value = ''
$ws.Cells.Item(238, 11).Value = 'Adverb / ਕਿਰਿਆ ਵਿਸੇਸ਼ਣ'
$ws.Cells.Item(238, 12).Value = 1
$ws.Cells.Item(238, 13).Value = 6
$ws.Cells.Item(238, 14).Value = 54730
$ws.Cells.Item(238, 15).Value = 1
$ws.Cells.Item(238, 16).Value = 1
$ws.Cells.Item(238, 17).Value = 3818
$ws.Cells.Item(238, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(238, 19).Value = ''
$ws.Cells.Item(238, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(238, 21).Value = ''
$ws.Cells.Item(238, 22).Value = ''
$ws.Cells.Item(238, 23).Value = ''
$ws.Cells.Item(238, 24).Value = ''
$ws.Cells.Item(238, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(238, 26).Value = 'Simple'
$ws.Cells.Item(238, 27).Value = 1376
$ws.Cells.Item(238, 28).Value = 0
$ws.Cells.Item(238, 29).Value = 0

# Row 239
$ws.Cells.Item(239, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(239, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(239, 3).Value = ''
$ws.Cells.Item(239, 4).Value = 'ਅਉਗਨੁ'
$ws.Cells.Item(239, 5).Value = 'No meanings found for ਅਉਗਨੁ'
$ws.Cells.Item(239, 6).Value = 'ੁ'
$ws.Cells.Item(239, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(239, 8).Value = 'Accusative ਕਰਮ'
$ws.Cells.Item(239, 9).Value = 'Masculine / ਪੁਲਿੰਗ'
$ws.Cells.Item(239, 10).Value = 'ਮੁਕਤਾ Ending'
$ws.Cells.Item(239, 11).Value = 'Noun / ਨਾਂਵ'
$ws.Cells.Item(239, 12).Value = 1
$ws.Cells.Item(239, 13).Value = 7
$ws.Cells.Item(239, 14).Value = 54730
$ws.Cells.Item(239, 15).Value = 1
$ws.Cells.Item(239, 16).Value = 1
$ws.Cells.Item(239, 17).Value = 3818
$ws.Cells.Item(239, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(239, 19).Value = ''
$ws.Cells.Item(239, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(239, 21).Value = ''
$ws.Cells.Item(239, 22).Value = ''
$ws.Cells.Item(239, 23).Value = ''
$ws.Cells.Item(239, 24).Value = ''
$ws.Cells.Item(239, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(239, 26).Value = 'Simple'
$ws.Cells.Item(239, 27).Value = 1376
$ws.Cells.Item(239, 28).Value = 0
$ws.Cells.Item(239, 29).Value = 0

# Row 240
$ws.Cells.Item(240, 1).Value = 'ਕਬੀਰ ਮਨੁ ਜਾਨੈ ਸਭ ਬਾਤ ਜਾਨਤ ਹੀ ਅਉਗਨੁ ਕਰੈ ॥'
$ws.Cells.Item(240, 2).Value = 'O Kabir! The mind knows all matter; even while knowing, one indeed commits the fault.'
$ws.Cells.Item(240, 3).Value = ''
$ws.Cells.Item(240, 4).Value = 'ਕਰੈ'
$ws.Cells.Item(240, 5).Value = 'ਕਰਦਾ ਹੈ| ਕਰਦਾ ਹੈ {ਇਕ-ਵਚਨ}'
$ws.Cells.Item(240, 6).Value = 'ਕਰੈ'
$ws.Cells.Item(240, 7).Value = 'Singular / ਇਕ'
$ws.Cells.Item(240, 8).Value = 'Present ਵਰਤਮਾਨ'
$ws.Cells.Item(240, 9).Value = 'Trans / ਨਪੁਂਸਕ'
$ws.Cells.Item(240, 10).Value = '3rd Person / ਅਨਯ ਪੁਰਖ'
$ws.Cells.Item(240, 11).Value = 'Verb / ਕਿਰਿਆ'
$ws.Cells.Item(240, 12).Value = 1
$ws.Cells.Item(240, 13).Value = 8
$ws.Cells.Item(240, 14).Value = 54730
$ws.Cells.Item(240, 15).Value = 1
$ws.Cells.Item(240, 16).Value = 1
$ws.Cells.Item(240, 17).Value = 3818
$ws.Cells.Item(240, 18).Value = 'ਸਲੋਕ ਭਗਤ ਕਬੀਰ ਜੀਉ ਕੇ '
$ws.Cells.Item(240, 19).Value = ''
$ws.Cells.Item(240, 20).Value = 'ਕਬੀਰ ਜੀ'
$ws.Cells.Item(240, 21).Value = ''
$ws.Cells.Item(240, 22).Value = ''
$ws.Cells.Item(240, 23).Value = ''
$ws.Cells.Item(240, 24).Value = ''
$ws.Cells.Item(240, 25).Value = 'ਸ਼ਲੋਕ'
$ws.Cells.Item(240, 26).Value = 'Simple'
$ws.Cells.Item(240, 27).Value = 1376
$ws.Cells.Item(240, 28).Value = 0
$ws.Cells.Item(240, 29).Value = 0
